# Populate rows 2-7 of Sheet1 per the JPM upcoming index events update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "AP"
$ws.Range("B2").Value = "4763 TT"
$ws.Range("C2").Value = "Jinan Acetate Chemical"
$ws.Range("D2").Value = "2025-07-07"
$ws.Range("E2").Value = "C"
$ws.Range("F2").Value = "FTSE EM Small Cap"
$ws.Range("G2").Value = "S Inc & F Inc"
$ws.Range("I2").Value = 2.844782917955921
$ws.Range("J2").Value = 0.09346468045032229
$ws.Range("K2").Value = 0.03019303399932429
$ws.Range("L2").Value = 2.844782917955921
$ws.Range("M2").Value = 0.09346468045032229
$ws.Range("N2").Value = 0.03019303399932429
$ws.Range("O2").Value = "2025-06-30"
$ws.Range("P2").Value = "Review"
$ws.Range("Q2").Value = "2025-07-01"

# Row 3
$ws.Range("A3").Value = "AP"
$ws.Range("B3").Value = "2888 TT"
$ws.Range("C3").Value = "Shin Kong Finl Hldgs"
$ws.Range("D3").Value = "2025-07-11"
$ws.Range("E3").Value = "E"
$ws.Range("F3").Value = "MSCI EM"
$ws.Range("G3").Value = "Delete"
$ws.Range("H3").Value = -0.0006745255708884523
$ws.Range("I3").Value = -31.54950164236167
$ws.Range("J3").Value = -81.61548344514789
$ws.Range("K3").Value = -1.304883791368007
$ws.Range("L3").Value = -31.54950164236167
$ws.Range("M3").Value = -81.61548344514789
$ws.Range("N3").Value = -1.304883791368007
$ws.Range("O3").Value = "2025-06-30"
$ws.Range("P3").Value = "Acquisition"
$ws.Range("Q3").Value = "2025-07-01"

# Row 4
$ws.Range("A4").Value = "AP"
$ws.Range("B4").Value = "2887 TT"
$ws.Range("C4").Value = "Taishin Financial Hldg"
$ws.Range("D4").Value = "2025-07-11"
$ws.Range("E4").Value = "E"
$ws.Range("F4").Value = "MSCI EM"
$ws.Range("G4").Value = "S Inc & F Dec"
$ws.Range("H4").Value = 0.0005781539159491507
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = "2025-06-30"
$ws.Range("P4").Value = "Acquisition"
$ws.Range("Q4").Value = "2025-07-01"

# Row 5
$ws.Range("A5").Value = "EMEA"
$ws.Range("B5").Value = "TBD"
$ws.Range("C5").Value = "Agility (Detached 2)"
$ws.Range("D5").Value = "2025-07-14"
$ws.Range("E5").Value = "E"
$ws.Range("F5").Value = "MSCI EM Small Cap"
$ws.Range("G5").Value = "Delete"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = "2025-06-23"
$ws.Range("P5").Value = "Spin-Off"
$ws.Range("Q5").Value = "2025-07-01"

# Row 6
$ws.Range("A6").Value = "EMEA"
$ws.Range("B6").Value = "AGILITY UH"
$ws.Range("C6").Value = "Agility Global"
$ws.Range("D6").Value = "2025-07-14"
$ws.Range("E6").Value = "E"
$ws.Range("F6").Value = "MSCI EM Small Cap"
$ws.Range("G6").Value = "F Inc"
$ws.Range("H6").Value = 0.0002387141887155395
$ws.Range("I6").Value = -0.2087005541376096
$ws.Range("J6").Value = -0.6440428949239859
$ws.Range("K6").Value = -0.06636583118172114
$ws.Range("L6").Value = -0.2087005541376096
$ws.Range("M6").Value = -0.6440428949239859
$ws.Range("N6").Value = -0.06636583118172114
$ws.Range("O6").Value = "2025-06-23"
$ws.Range("P6").Value = "Spin-Off"
$ws.Range("Q6").Value = "2025-07-01"

# Row 7
$ws.Range("A7").Value = "AP"
$ws.Range("B7").Value = "9787 JP"
$ws.Range("C7").Value = "Aeon Delight Co Ltd"
$ws.Range("D7").Value = "2025-07-16"
$ws.Range("E7").Value = "C"
$ws.Range("F7").Value = "FTSE DM Small Cap"
$ws.Range("G7").Value = "Delete"
$ws.Range("I7").Value = -7.102833794188818
$ws.Range("J7").Value = -0.1920675254749968
$ws.Range("K7").Value = -2.436014020863679
$ws.Range("L7").Value = -7.102833794188818
$ws.Range("M7").Value = -0.1920675254749968
$ws.Range("N7").Value = -2.436014020863679
$ws.Range("O7").Value = "2025-06-25"
$ws.Range("P7").Value = "Delisting"
$ws.Range("Q7").Value = "2025-07-01"

# The source data stores these flag cells as the literal text "TRUE" (not a
# Boolean) -- plain `.Value = "TRUE"` gets auto-coerced by Excel into a real
# Boolean, so a leading apostrophe (quote-prefix) forces text entry; resetting
# the style back to Normal afterwards drops the quote-prefix cell formatting
# so only the literal string value remains on the cell.
$trueTextCells = @("R2", "T2", "W3", "R4", "U4", "W5", "X5", "T6", "X6", "W7")
foreach ($ref in $trueTextCells) {
    $ws.Range($ref).Value = "'TRUE"
    $ws.Range($ref).Style = "Normal"
}
